# Apply "Project - final iteration" edits to the technology sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update part identity & first operation ---
$ws.Range("A2").Value = "DA08-00.0100.02"
$ws.Range("B2").Value = "Корпус статора"
$ws.Range("C2").Value = "S17A.6113"
$ws.Range("H2").Value = "OP10"
$ws.Range("I2").Value = "Ленточнопильная"
$ws.Range("J2").Value = "FMB"
$ws.Range("L2").Value = "2,0"

# --- Row 3: second operation ---
$ws.Range("H3").Value = "OP40"
$ws.Range("I3").Value = "Зачистная"
$ws.Range("J3").Value = "Bomar Single"

# --- Row 4: third operation ---
$ws.Range("H4").Value = "OP350"
$ws.Range("I4").Value = "Очистка и консервахция"
$ws.Range("J4").Value = "ELMA Xхtra line AM"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "0"
$ws.Range("K4").ClearFormats()
$ws.Range("L4").Value = "15,0"
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "12"
$ws.Range("M4").ClearFormats()

# --- Row 5: fourth operation ---
$ws.Range("H5").Value = "OP360"
$ws.Range("I5").Value = "Укладка в складскую тару"
$ws.Range("J5").Value = "Packager's WorkPlace"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "0"
$ws.Range("K5").ClearFormats()
$ws.Range("L5").Value = "0,25"
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "1"
$ws.Range("M5").ClearFormats()

# --- Remove rows 6-11 (no longer part of the data) ---
$ws.Range("A6:M11").EntireRow.Delete()
